# PowerShell-style Excel COM-interop script
# Applies the "output generated at 456a3b4" data refresh:
#  - updates "想去人数" (F) / occasional "最低票价" (G) counters on 展览 & 全部类型
#  - marks the 花样年华 concert (row 3) as cancelled/sold-out on 演出 & 全部类型
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 3150
$wsExhibit.Range("F5").Value = 2170
$wsExhibit.Range("F7").Value = 309
$wsExhibit.Range("F8").Value = 929
$wsExhibit.Range("F9").Value = 1000
$wsExhibit.Range("F10").Value = 237
$wsExhibit.Range("F11").Value = 457
$wsExhibit.Range("F12").Value = 1154
$wsExhibit.Range("F16").Value = 7743
$wsExhibit.Range("F17").Value = 338
$wsExhibit.Range("F18").Value = 2456
$wsExhibit.Range("F19").Value = 212
$wsExhibit.Range("F20").Value = 226
$wsExhibit.Range("F22").Value = 458
$wsExhibit.Range("F28").Value = 1650
$wsExhibit.Range("G28").Value = 58.5
$wsExhibit.Range("F29").Value = 240
$wsExhibit.Range("F30").Value = 1167
$wsExhibit.Range("F33").Value = 37
$wsExhibit.Range("F34").Value = 165
$wsExhibit.Range("F35").Value = 270
$wsExhibit.Range("F36").Value = 38
$wsExhibit.Range("F37").Value = 163
$wsExhibit.Range("F38").Value = 336
$wsExhibit.Range("F40").Value = 215

# --- Sheet "演出" (Performances) - cancel 花样年华 concert ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("C3").Value = "浙江·花样年华·一生必听的经典电影名曲——《花样年华》《泰坦尼克号》《爱乐之城》（取消）"
$wsShow.Range("G3").Value = "不可售"

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("C3").Value = "浙江·花样年华·一生必听的经典电影名曲——《花样年华》《泰坦尼克号》《爱乐之城》（取消）"
$wsAll.Range("G3").Value = "不可售"
$wsAll.Range("F5").Value = 3150
$wsAll.Range("F7").Value = 2170
$wsAll.Range("F9").Value = 309
$wsAll.Range("F10").Value = 929
$wsAll.Range("F12").Value = 1000
$wsAll.Range("F13").Value = 237
$wsAll.Range("F14").Value = 457
$wsAll.Range("F15").Value = 1154
$wsAll.Range("F19").Value = 7743
$wsAll.Range("F20").Value = 338
$wsAll.Range("F21").Value = 2456
$wsAll.Range("F23").Value = 212
$wsAll.Range("F24").Value = 226
$wsAll.Range("F26").Value = 458
$wsAll.Range("F32").Value = 1650
$wsAll.Range("G32").Value = 58.5
$wsAll.Range("F33").Value = 240
$wsAll.Range("F34").Value = 1167
$wsAll.Range("F37").Value = 37
$wsAll.Range("F38").Value = 165
$wsAll.Range("F39").Value = 270
$wsAll.Range("F40").Value = 38
$wsAll.Range("F41").Value = 163
$wsAll.Range("F42").Value = 336
$wsAll.Range("F47").Value = 215
